$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC003 "Expected Result" cell (C4): append a third accepted error message
# ("Không được để trống số điện thoại") to the existing text.
$ws.Range("C4").Value = "The system shows error message  `n""SĐT phải là 10 số"" or `n""SĐT không hợp lệ""or`n""Không được để trống số điện thoại"""

# TC001 "Method" cell (G2) switched from Auto to Manual
$ws.Range("G2").Value = "Manual"

# Column D (Test Suite) widened slightly to fit its contents
$ws.Columns.Item(4).ColumnWidth = 30.5

# Row 4 grows taller to accommodate the now-longer wrapped text in C4
$ws.Rows.Item(4).RowHeight = 69

# Leave the selection where the editor last clicked, D3
$ws.Range("D3").Select()
